# Actualizacion final de lista de tareas Rev2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Status column (F) for rows 6,8,9,10,11 moves from "Por iniciar" to the
# new status "Terminado". Row 7's status text is renamed "Por iniciar" -> "En proceso".
$ws.Range("F7").Value = "En proceso"
$ws.Range("F6").Value = "Terminado"
$ws.Range("F8").Value = "Terminado"
$ws.Range("F9").Value = "Terminado"
$ws.Range("F10").Value = "Terminado"
$ws.Range("F11").Value = "Terminado"

# "Horas reales" / remaining hours entries (AZ column) for rows 6-10.
$ws.Range("AZ6").Value = 4
$ws.Range("AZ7").Value = 2.5
$ws.Range("AZ8").Value = 4
$ws.Range("AZ9").Value = 2
$ws.Range("AZ10").Value = 4

# View settings: zoom back to 100%, and selection moved to F7.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("F7").Select()
